$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Jumlah" (Total) rows: sum the alternating data rows above each total row.
# Match the "Jumlah" label cell's formatting (C29 / C30) across the whole
# total row before writing the formulas, same as the source workbook.

# Row 29 totals the "Baru" (new) rows: 17, 19, 21, 23, 25, 27
$ws.Range("C29").Copy()
$ws.Range("D29:AI29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D29:AI29").Formula = "=SUM(D17,D19,D21,D23,D25,D27)"

# Row 30 totals the "Ulangan" (repeat) rows: 18, 20, 22, 24, 26, 28
$ws.Range("C30").Copy()
$ws.Range("D30:AI30").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D30:AI30").Formula = "=SUM(D18,D20,D22,D24,D26,D28)"

$excel.CutCopyMode = $false

# Restore the view state: scrolled down a bit with D30:AI30 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("D30:AI30").Select()
